$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete old row 65 (Trend Bluetooth Lautsprecher E100 Gruen) - shifts rows up
$ws.Rows.Item(65).Delete()

# Step 2: update timestamp column (O) for all data rows to the new value
$newTimestamp = "2022-07-28 20:59:25"
for ($r = 2; $r -le 86; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Step 3: rewrite rows whose content moved/changed due to reordering
# Row 15: '4949707' - 'Varta Knopfzellen CR2025 2 Stück'
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "4949707"
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "Varta Knopfzellen CR2025 2 Stück"
$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-knopfzellen-cr2025-2-stueck/p/4949707"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "2ST"
$ws.Cells.Item(15, 5).ClearContents()
$ws.Cells.Item(15, 6).NumberFormat = "General"
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "Varta"
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value = "8.95"
$ws.Cells.Item(15, 9).NumberFormat = "@"
$ws.Cells.Item(15, 9).Value = "4.48/1ST"
$ws.Cells.Item(15, 10).NumberFormat = "@"
$ws.Cells.Item(15, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(15, 11).NumberFormat = "@"
$ws.Cells.Item(15, 11).Value = "4.48"
$ws.Cells.Item(15, 12).NumberFormat = "@"
$ws.Cells.Item(15, 12).Value = "1ST"
$ws.Cells.Item(15, 13).NumberFormat = "@"
$ws.Cells.Item(15, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(15, 14).NumberFormat = "@"
$ws.Cells.Item(15, 14).Value = "Varta Knopfzellen CR2025 2 Stück 8.95 Schweizer Franken"
$ws.Cells.Item(15, 15).NumberFormat = "@"
$ws.Cells.Item(15, 15).Value = "2022-07-28 20:59:25"

# Row 16: '6153846' - 'Duracell Knopfzelle CR2032 3V 2 Stück'
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "6153846"
$ws.Cells.Item(16, 2).NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = "Duracell Knopfzelle CR2032 3V 2 Stück"
$ws.Cells.Item(16, 3).NumberFormat = "@"
$ws.Cells.Item(16, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-knopfzelle-cr2032-3v-2-stueck/p/6153846"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2ST"
$ws.Cells.Item(16, 5).NumberFormat = "General"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).NumberFormat = "General"
$ws.Cells.Item(16, 6).Value = 5
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "Duracell"
$ws.Cells.Item(16, 8).NumberFormat = "@"
$ws.Cells.Item(16, 8).Value = "9.95"
$ws.Cells.Item(16, 9).NumberFormat = "@"
$ws.Cells.Item(16, 9).Value = "4.98/1ST"
$ws.Cells.Item(16, 10).NumberFormat = "@"
$ws.Cells.Item(16, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(16, 11).NumberFormat = "@"
$ws.Cells.Item(16, 11).Value = "4.98"
$ws.Cells.Item(16, 12).NumberFormat = "@"
$ws.Cells.Item(16, 12).Value = "1ST"
$ws.Cells.Item(16, 13).NumberFormat = "@"
$ws.Cells.Item(16, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(16, 14).NumberFormat = "@"
$ws.Cells.Item(16, 14).Value = "Duracell Knopfzelle CR2032 3V 2 Stück 9.95 Schweizer Franken"
$ws.Cells.Item(16, 15).NumberFormat = "@"
$ws.Cells.Item(16, 15).Value = "2022-07-28 20:59:25"

# Row 22: '6761133' - 'Duracell Batterien PLUS C/LR14 2 Stück'
$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = "6761133"
$ws.Cells.Item(22, 2).NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = "Duracell Batterien PLUS C/LR14 2 Stück"
$ws.Cells.Item(22, 3).NumberFormat = "@"
$ws.Cells.Item(22, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterien-plus-clr14-2-stueck/p/6761133"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "2ST"
$ws.Cells.Item(22, 5).ClearContents()
$ws.Cells.Item(22, 6).NumberFormat = "General"
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "Duracell"
$ws.Cells.Item(22, 8).NumberFormat = "@"
$ws.Cells.Item(22, 8).Value = "9.95"
$ws.Cells.Item(22, 9).NumberFormat = "@"
$ws.Cells.Item(22, 9).Value = "4.98/1ST"
$ws.Cells.Item(22, 10).NumberFormat = "@"
$ws.Cells.Item(22, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(22, 11).NumberFormat = "@"
$ws.Cells.Item(22, 11).Value = "4.98"
$ws.Cells.Item(22, 12).NumberFormat = "@"
$ws.Cells.Item(22, 12).Value = "1ST"
$ws.Cells.Item(22, 13).NumberFormat = "@"
$ws.Cells.Item(22, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(22, 14).NumberFormat = "@"
$ws.Cells.Item(22, 14).Value = "Duracell Batterien PLUS C/LR14 2 Stück 9.95 Schweizer Franken"
$ws.Cells.Item(22, 15).NumberFormat = "@"
$ws.Cells.Item(22, 15).Value = "2022-07-28 20:59:25"

# Row 23: '6761134' - 'Duracell Batterien PLUS D/LR20 2 Stück'
$ws.Cells.Item(23, 1).NumberFormat = "@"
$ws.Cells.Item(23, 1).Value = "6761134"
$ws.Cells.Item(23, 2).NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = "Duracell Batterien PLUS D/LR20 2 Stück"
$ws.Cells.Item(23, 3).NumberFormat = "@"
$ws.Cells.Item(23, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterien-plus-dlr20-2-stueck/p/6761134"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2ST"
$ws.Cells.Item(23, 5).ClearContents()
$ws.Cells.Item(23, 6).NumberFormat = "General"
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "Duracell"
$ws.Cells.Item(23, 8).NumberFormat = "@"
$ws.Cells.Item(23, 8).Value = "9.95"
$ws.Cells.Item(23, 9).NumberFormat = "@"
$ws.Cells.Item(23, 9).Value = "4.98/1ST"
$ws.Cells.Item(23, 10).NumberFormat = "@"
$ws.Cells.Item(23, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(23, 11).NumberFormat = "@"
$ws.Cells.Item(23, 11).Value = "4.98"
$ws.Cells.Item(23, 12).NumberFormat = "@"
$ws.Cells.Item(23, 12).Value = "1ST"
$ws.Cells.Item(23, 13).NumberFormat = "@"
$ws.Cells.Item(23, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(23, 14).NumberFormat = "@"
$ws.Cells.Item(23, 14).Value = "Duracell Batterien PLUS D/LR20 2 Stück 9.95 Schweizer Franken"
$ws.Cells.Item(23, 15).NumberFormat = "@"
$ws.Cells.Item(23, 15).Value = "2022-07-28 20:59:25"

# Row 24: '4014527' - 'Varta Longlife Batterien AA/LR6 10 Stück'
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "4014527"
$ws.Cells.Item(24, 2).NumberFormat = "@"
$ws.Cells.Item(24, 2).Value = "Varta Longlife Batterien AA/LR6 10 Stück"
$ws.Cells.Item(24, 3).NumberFormat = "@"
$ws.Cells.Item(24, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-batterien-aalr6-10-stueck/p/4014527"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "10ST"
$ws.Cells.Item(24, 5).NumberFormat = "General"
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).NumberFormat = "General"
$ws.Cells.Item(24, 6).Value = 3
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "Varta"
$ws.Cells.Item(24, 8).NumberFormat = "@"
$ws.Cells.Item(24, 8).Value = "15.95"
$ws.Cells.Item(24, 9).NumberFormat = "@"
$ws.Cells.Item(24, 9).Value = "1.60/1ST"
$ws.Cells.Item(24, 10).NumberFormat = "@"
$ws.Cells.Item(24, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(24, 11).NumberFormat = "@"
$ws.Cells.Item(24, 11).Value = "1.60"
$ws.Cells.Item(24, 12).NumberFormat = "@"
$ws.Cells.Item(24, 12).Value = "1ST"
$ws.Cells.Item(24, 13).NumberFormat = "@"
$ws.Cells.Item(24, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Cells.Item(24, 14).NumberFormat = "@"
$ws.Cells.Item(24, 14).Value = "Varta Longlife Batterien AA/LR6 10 Stück 15.95 Schweizer Franken"
$ws.Cells.Item(24, 15).NumberFormat = "@"
$ws.Cells.Item(24, 15).Value = "2022-07-28 20:59:25"

# Row 41: '3494063' - 'Varta Longlife Power C 2er Bli'
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = "3494063"
$ws.Cells.Item(41, 2).NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = "Varta Longlife Power C 2er Bli"
$ws.Cells.Item(41, 3).NumberFormat = "@"
$ws.Cells.Item(41, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-c-2er-bli/p/3494063"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2ST"
$ws.Cells.Item(41, 5).ClearContents()
$ws.Cells.Item(41, 6).NumberFormat = "General"
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).NumberFormat = "@"
$ws.Cells.Item(41, 7).Value = "Varta"
$ws.Cells.Item(41, 8).NumberFormat = "@"
$ws.Cells.Item(41, 8).Value = "7.95"
$ws.Cells.Item(41, 9).NumberFormat = "@"
$ws.Cells.Item(41, 9).Value = "3.98/1ST"
$ws.Cells.Item(41, 10).NumberFormat = "@"
$ws.Cells.Item(41, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(41, 11).NumberFormat = "@"
$ws.Cells.Item(41, 11).Value = "3.98"
$ws.Cells.Item(41, 12).NumberFormat = "@"
$ws.Cells.Item(41, 12).Value = "1ST"
$ws.Cells.Item(41, 13).NumberFormat = "@"
$ws.Cells.Item(41, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(41, 14).NumberFormat = "@"
$ws.Cells.Item(41, 14).Value = "Varta Longlife Power C 2er Bli 7.95 Schweizer Franken"
$ws.Cells.Item(41, 15).NumberFormat = "@"
$ws.Cells.Item(41, 15).Value = "2022-07-28 20:59:25"

# Row 42: '4358323' - 'Rayovac Hörgerätebatterien 312 6 Stück'
$ws.Cells.Item(42, 1).NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = "4358323"
$ws.Cells.Item(42, 2).NumberFormat = "@"
$ws.Cells.Item(42, 2).Value = "Rayovac Hörgerätebatterien 312 6 Stück"
$ws.Cells.Item(42, 3).NumberFormat = "@"
$ws.Cells.Item(42, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/rayovac-hoergeraetebatterien-312-6-stueck/p/4358323"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "6ST"
$ws.Cells.Item(42, 5).NumberFormat = "General"
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(42, 6).NumberFormat = "General"
$ws.Cells.Item(42, 6).Value = 4
$ws.Cells.Item(42, 7).NumberFormat = "@"
$ws.Cells.Item(42, 7).Value = "Rayovac"
$ws.Cells.Item(42, 8).NumberFormat = "@"
$ws.Cells.Item(42, 8).Value = "9.95"
$ws.Cells.Item(42, 9).NumberFormat = "@"
$ws.Cells.Item(42, 9).Value = "1.66/1ST"
$ws.Cells.Item(42, 10).NumberFormat = "@"
$ws.Cells.Item(42, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(42, 11).NumberFormat = "@"
$ws.Cells.Item(42, 11).Value = "1.66"
$ws.Cells.Item(42, 12).NumberFormat = "@"
$ws.Cells.Item(42, 12).Value = "1ST"
$ws.Cells.Item(42, 13).NumberFormat = "@"
$ws.Cells.Item(42, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(42, 14).NumberFormat = "@"
$ws.Cells.Item(42, 14).Value = "Rayovac Hörgerätebatterien 312 6 Stück 9.95 Schweizer Franken"
$ws.Cells.Item(42, 15).NumberFormat = "@"
$ws.Cells.Item(42, 15).Value = "2022-07-28 20:59:25"

# Row 66: '6872591' - 'Bosch Küchenmaschine MUM58243 1000W'
$ws.Cells.Item(66, 1).NumberFormat = "@"
$ws.Cells.Item(66, 1).Value = "6872591"
$ws.Cells.Item(66, 2).NumberFormat = "@"
$ws.Cells.Item(66, 2).Value = "Bosch Küchenmaschine MUM58243 1000W"
$ws.Cells.Item(66, 3).NumberFormat = "@"
$ws.Cells.Item(66, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/bosch-kuechenmaschine-mum58243-1000w/p/6872591"
$ws.Cells.Item(66, 4).ClearContents()
$ws.Cells.Item(66, 5).ClearContents()
$ws.Cells.Item(66, 6).NumberFormat = "General"
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).NumberFormat = "@"
$ws.Cells.Item(66, 7).Value = "Bosch"
$ws.Cells.Item(66, 8).NumberFormat = "@"
$ws.Cells.Item(66, 8).Value = "239.50"
$ws.Cells.Item(66, 9).ClearContents()
$ws.Cells.Item(66, 10).ClearContents()
$ws.Cells.Item(66, 11).ClearContents()
$ws.Cells.Item(66, 12).ClearContents()
$ws.Cells.Item(66, 13).NumberFormat = "@"
$ws.Cells.Item(66, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Cells.Item(66, 14).NumberFormat = "@"
$ws.Cells.Item(66, 14).Value = "Bosch Küchenmaschine MUM58243 1000W 50% Aktion 239.50 Schweizer Franken statt 479.00 Schweizer Franken"
$ws.Cells.Item(66, 15).NumberFormat = "@"
$ws.Cells.Item(66, 15).Value = "2022-07-28 20:59:25"

# Row 67: '6986541' - 'Duracell Batterie (CR2032, 4 Stück)'
$ws.Cells.Item(67, 1).NumberFormat = "@"
$ws.Cells.Item(67, 1).Value = "6986541"
$ws.Cells.Item(67, 2).NumberFormat = "@"
$ws.Cells.Item(67, 2).Value = "Duracell Batterie (CR2032, 4 Stück)"
$ws.Cells.Item(67, 3).NumberFormat = "@"
$ws.Cells.Item(67, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterie-cr2032-4-stueck/p/6986541"
$ws.Cells.Item(67, 4).NumberFormat = "@"
$ws.Cells.Item(67, 4).Value = "4ST"
$ws.Cells.Item(67, 5).ClearContents()
$ws.Cells.Item(67, 6).NumberFormat = "General"
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).NumberFormat = "@"
$ws.Cells.Item(67, 7).Value = "Duracell"
$ws.Cells.Item(67, 8).NumberFormat = "@"
$ws.Cells.Item(67, 8).Value = "8.95"
$ws.Cells.Item(67, 9).NumberFormat = "@"
$ws.Cells.Item(67, 9).Value = "2.24/1ST"
$ws.Cells.Item(67, 10).NumberFormat = "@"
$ws.Cells.Item(67, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(67, 11).NumberFormat = "@"
$ws.Cells.Item(67, 11).Value = "2.24"
$ws.Cells.Item(67, 12).NumberFormat = "@"
$ws.Cells.Item(67, 12).Value = "1ST"
$ws.Cells.Item(67, 13).NumberFormat = "@"
$ws.Cells.Item(67, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(67, 14).NumberFormat = "@"
$ws.Cells.Item(67, 14).Value = "Duracell Batterie (CR2032, 4 Stück) 43% Aktion 8.95 Schweizer Franken statt 15.90 Schweizer Franken"
$ws.Cells.Item(67, 15).NumberFormat = "@"
$ws.Cells.Item(67, 15).Value = "2022-07-28 20:59:25"

# Row 68: '6867383' - 'Krups Kaffeevollautomat EA815B'
$ws.Cells.Item(68, 1).NumberFormat = "@"
$ws.Cells.Item(68, 1).Value = "6867383"
$ws.Cells.Item(68, 2).NumberFormat = "@"
$ws.Cells.Item(68, 2).Value = "Krups Kaffeevollautomat EA815B"
$ws.Cells.Item(68, 3).NumberFormat = "@"
$ws.Cells.Item(68, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/krups-kaffeevollautomat-ea815b/p/6867383"
$ws.Cells.Item(68, 4).ClearContents()
$ws.Cells.Item(68, 5).ClearContents()
$ws.Cells.Item(68, 6).NumberFormat = "General"
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).NumberFormat = "@"
$ws.Cells.Item(68, 7).Value = "Krups"
$ws.Cells.Item(68, 8).NumberFormat = "@"
$ws.Cells.Item(68, 8).Value = "349.50"
$ws.Cells.Item(68, 9).ClearContents()
$ws.Cells.Item(68, 10).ClearContents()
$ws.Cells.Item(68, 11).ClearContents()
$ws.Cells.Item(68, 12).ClearContents()
$ws.Cells.Item(68, 13).NumberFormat = "@"
$ws.Cells.Item(68, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Cells.Item(68, 14).NumberFormat = "@"
$ws.Cells.Item(68, 14).Value = "Krups Kaffeevollautomat EA815B 50% Aktion 349.50 Schweizer Franken statt 699.00 Schweizer Franken"
$ws.Cells.Item(68, 15).NumberFormat = "@"
$ws.Cells.Item(68, 15).Value = "2022-07-28 20:59:25"

# Row 69: '6735643' - 'LED 31V Anschlussset Transf.+Verl.kabel'
$ws.Cells.Item(69, 1).NumberFormat = "@"
$ws.Cells.Item(69, 1).Value = "6735643"
$ws.Cells.Item(69, 2).NumberFormat = "@"
$ws.Cells.Item(69, 2).Value = "LED 31V Anschlussset Transf.+Verl.kabel"
$ws.Cells.Item(69, 3).NumberFormat = "@"
$ws.Cells.Item(69, 3).Value = "/de/haushalt-tier/haushalt-kueche/uebrige-haushaltsartikel/led-31v-anschlussset-transfverlkabel/p/6735643"
$ws.Cells.Item(69, 4).ClearContents()
$ws.Cells.Item(69, 5).NumberFormat = "General"
$ws.Cells.Item(69, 5).Value = 1
$ws.Cells.Item(69, 6).NumberFormat = "General"
$ws.Cells.Item(69, 6).Value = 5
$ws.Cells.Item(69, 7).NumberFormat = "@"
$ws.Cells.Item(69, 7).Value = "Coop"
$ws.Cells.Item(69, 8).NumberFormat = "@"
$ws.Cells.Item(69, 8).Value = "9.95"
$ws.Cells.Item(69, 9).ClearContents()
$ws.Cells.Item(69, 10).ClearContents()
$ws.Cells.Item(69, 11).ClearContents()
$ws.Cells.Item(69, 12).ClearContents()
$ws.Cells.Item(69, 13).NumberFormat = "@"
$ws.Cells.Item(69, 13).Value = "['haushalt-tier', 'haushalt-kueche', 'uebrige-haushaltsartikel']"
$ws.Cells.Item(69, 14).NumberFormat = "@"
$ws.Cells.Item(69, 14).Value = "LED 31V Anschlussset Transf.+Verl.kabel 50% Aktion 9.95 Schweizer Franken statt 19.95 Schweizer Franken"
$ws.Cells.Item(69, 15).NumberFormat = "@"
$ws.Cells.Item(69, 15).Value = "2022-07-28 20:59:25"

# Row 70: '5882124' - 'Philips Avent Audio Monitors DECT-Babyphone'
$ws.Cells.Item(70, 1).NumberFormat = "@"
$ws.Cells.Item(70, 1).Value = "5882124"
$ws.Cells.Item(70, 2).NumberFormat = "@"
$ws.Cells.Item(70, 2).Value = "Philips Avent Audio Monitors DECT-Babyphone"
$ws.Cells.Item(70, 3).NumberFormat = "@"
$ws.Cells.Item(70, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/philips-avent-audio-monitors-dect-babyphone/p/5882124"
$ws.Cells.Item(70, 4).ClearContents()
$ws.Cells.Item(70, 5).ClearContents()
$ws.Cells.Item(70, 6).NumberFormat = "General"
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).NumberFormat = "@"
$ws.Cells.Item(70, 7).Value = "Avent"
$ws.Cells.Item(70, 8).NumberFormat = "@"
$ws.Cells.Item(70, 8).Value = "99.90"
$ws.Cells.Item(70, 9).ClearContents()
$ws.Cells.Item(70, 10).ClearContents()
$ws.Cells.Item(70, 11).ClearContents()
$ws.Cells.Item(70, 12).ClearContents()
$ws.Cells.Item(70, 13).NumberFormat = "@"
$ws.Cells.Item(70, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete']"
$ws.Cells.Item(70, 14).NumberFormat = "@"
$ws.Cells.Item(70, 14).Value = "Philips Avent Audio Monitors DECT-Babyphone 99.90 Schweizer Franken"
$ws.Cells.Item(70, 15).NumberFormat = "@"
$ws.Cells.Item(70, 15).Value = "2022-07-28 20:59:25"

# Row 71: '6425996' - 'satrap espresso XA Kolbenkaffeemaschine'
$ws.Cells.Item(71, 1).NumberFormat = "@"
$ws.Cells.Item(71, 1).Value = "6425996"
$ws.Cells.Item(71, 2).NumberFormat = "@"
$ws.Cells.Item(71, 2).Value = "satrap espresso XA Kolbenkaffeemaschine"
$ws.Cells.Item(71, 3).NumberFormat = "@"
$ws.Cells.Item(71, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-espresso-xa-kolbenkaffeemaschine/p/6425996"
$ws.Cells.Item(71, 4).ClearContents()
$ws.Cells.Item(71, 5).NumberFormat = "General"
$ws.Cells.Item(71, 5).Value = 7
$ws.Cells.Item(71, 6).NumberFormat = "General"
$ws.Cells.Item(71, 6).Value = 4
$ws.Cells.Item(71, 7).NumberFormat = "@"
$ws.Cells.Item(71, 7).Value = "satrap"
$ws.Cells.Item(71, 8).NumberFormat = "@"
$ws.Cells.Item(71, 8).Value = "249.00"
$ws.Cells.Item(71, 9).ClearContents()
$ws.Cells.Item(71, 10).ClearContents()
$ws.Cells.Item(71, 11).ClearContents()
$ws.Cells.Item(71, 12).ClearContents()
$ws.Cells.Item(71, 13).NumberFormat = "@"
$ws.Cells.Item(71, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Cells.Item(71, 14).NumberFormat = "@"
$ws.Cells.Item(71, 14).Value = "satrap espresso XA Kolbenkaffeemaschine 249.00 Schweizer Franken"
$ws.Cells.Item(71, 15).NumberFormat = "@"
$ws.Cells.Item(71, 15).Value = "2022-07-28 20:59:25"

# Row 72: '6125818' - 'satrap Mano XA Handmixer'
$ws.Cells.Item(72, 1).NumberFormat = "@"
$ws.Cells.Item(72, 1).Value = "6125818"
$ws.Cells.Item(72, 2).NumberFormat = "@"
$ws.Cells.Item(72, 2).Value = "satrap Mano XA Handmixer"
$ws.Cells.Item(72, 3).NumberFormat = "@"
$ws.Cells.Item(72, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-mano-xa-handmixer/p/6125818"
$ws.Cells.Item(72, 4).ClearContents()
$ws.Cells.Item(72, 5).ClearContents()
$ws.Cells.Item(72, 6).NumberFormat = "General"
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).NumberFormat = "@"
$ws.Cells.Item(72, 7).Value = "satrap"
$ws.Cells.Item(72, 8).NumberFormat = "@"
$ws.Cells.Item(72, 8).Value = "49.95"
$ws.Cells.Item(72, 9).ClearContents()
$ws.Cells.Item(72, 10).ClearContents()
$ws.Cells.Item(72, 11).ClearContents()
$ws.Cells.Item(72, 12).ClearContents()
$ws.Cells.Item(72, 13).NumberFormat = "@"
$ws.Cells.Item(72, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Cells.Item(72, 14).NumberFormat = "@"
$ws.Cells.Item(72, 14).Value = "satrap Mano XA Handmixer 49.95 Schweizer Franken"
$ws.Cells.Item(72, 15).NumberFormat = "@"
$ws.Cells.Item(72, 15).Value = "2022-07-28 20:59:25"

# Row 73: '5831402' - 'Satrap Mikrowelle Micro M2'
$ws.Cells.Item(73, 1).NumberFormat = "@"
$ws.Cells.Item(73, 1).Value = "5831402"
$ws.Cells.Item(73, 2).NumberFormat = "@"
$ws.Cells.Item(73, 2).Value = "Satrap Mikrowelle Micro M2"
$ws.Cells.Item(73, 3).NumberFormat = "@"
$ws.Cells.Item(73, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-mikrowelle-micro-m2/p/5831402"
$ws.Cells.Item(73, 4).ClearContents()
$ws.Cells.Item(73, 5).ClearContents()
$ws.Cells.Item(73, 6).NumberFormat = "General"
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).NumberFormat = "@"
$ws.Cells.Item(73, 7).Value = "satrap"
$ws.Cells.Item(73, 8).NumberFormat = "@"
$ws.Cells.Item(73, 8).Value = "49.95"
$ws.Cells.Item(73, 9).ClearContents()
$ws.Cells.Item(73, 10).ClearContents()
$ws.Cells.Item(73, 11).ClearContents()
$ws.Cells.Item(73, 12).ClearContents()
$ws.Cells.Item(73, 13).NumberFormat = "@"
$ws.Cells.Item(73, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Cells.Item(73, 14).NumberFormat = "@"
$ws.Cells.Item(73, 14).Value = "Satrap Mikrowelle Micro M2 50% Aktion 49.95 Schweizer Franken statt 99.90 Schweizer Franken"
$ws.Cells.Item(73, 15).NumberFormat = "@"
$ws.Cells.Item(73, 15).Value = "2022-07-28 20:59:25"

# Row 74: '5872164' - 'Satrap Tischventilator Venti 2'
$ws.Cells.Item(74, 1).NumberFormat = "@"
$ws.Cells.Item(74, 1).Value = "5872164"
$ws.Cells.Item(74, 2).NumberFormat = "@"
$ws.Cells.Item(74, 2).Value = "Satrap Tischventilator Venti 2"
$ws.Cells.Item(74, 3).NumberFormat = "@"
$ws.Cells.Item(74, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-tischventilator-venti-2/p/5872164"
$ws.Cells.Item(74, 4).ClearContents()
$ws.Cells.Item(74, 5).ClearContents()
$ws.Cells.Item(74, 6).NumberFormat = "General"
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).NumberFormat = "@"
$ws.Cells.Item(74, 7).Value = "satrap"
$ws.Cells.Item(74, 8).NumberFormat = "@"
$ws.Cells.Item(74, 8).Value = "34.95"
$ws.Cells.Item(74, 9).ClearContents()
$ws.Cells.Item(74, 10).ClearContents()
$ws.Cells.Item(74, 11).ClearContents()
$ws.Cells.Item(74, 12).ClearContents()
$ws.Cells.Item(74, 13).NumberFormat = "@"
$ws.Cells.Item(74, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Cells.Item(74, 14).NumberFormat = "@"
$ws.Cells.Item(74, 14).Value = "Satrap Tischventilator Venti 2 34.95 Schweizer Franken"
$ws.Cells.Item(74, 15).NumberFormat = "@"
$ws.Cells.Item(74, 15).Value = "2022-07-28 20:59:25"

# Row 75: '6822783' - 'Trend Car Charger QC'
$ws.Cells.Item(75, 1).NumberFormat = "@"
$ws.Cells.Item(75, 1).Value = "6822783"
$ws.Cells.Item(75, 2).NumberFormat = "@"
$ws.Cells.Item(75, 2).Value = "Trend Car Charger QC"
$ws.Cells.Item(75, 3).NumberFormat = "@"
$ws.Cells.Item(75, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-car-charger-qc/p/6822783"
$ws.Cells.Item(75, 4).ClearContents()
$ws.Cells.Item(75, 5).ClearContents()
$ws.Cells.Item(75, 6).NumberFormat = "General"
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).NumberFormat = "@"
$ws.Cells.Item(75, 7).Value = "Trend"
$ws.Cells.Item(75, 8).NumberFormat = "@"
$ws.Cells.Item(75, 8).Value = "19.95"
$ws.Cells.Item(75, 9).ClearContents()
$ws.Cells.Item(75, 10).ClearContents()
$ws.Cells.Item(75, 11).ClearContents()
$ws.Cells.Item(75, 12).ClearContents()
$ws.Cells.Item(75, 13).NumberFormat = "@"
$ws.Cells.Item(75, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Cells.Item(75, 14).NumberFormat = "@"
$ws.Cells.Item(75, 14).Value = "Trend Car Charger QC 19.95 Schweizer Franken"
$ws.Cells.Item(75, 15).NumberFormat = "@"
$ws.Cells.Item(75, 15).Value = "2022-07-28 20:59:25"

# Row 76: '6007537' - 'Trend USB-Stick 16 GB'
$ws.Cells.Item(76, 1).NumberFormat = "@"
$ws.Cells.Item(76, 1).Value = "6007537"
$ws.Cells.Item(76, 2).NumberFormat = "@"
$ws.Cells.Item(76, 2).Value = "Trend USB-Stick 16 GB"
$ws.Cells.Item(76, 3).NumberFormat = "@"
$ws.Cells.Item(76, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-16-gb/p/6007537"
$ws.Cells.Item(76, 4).ClearContents()
$ws.Cells.Item(76, 5).NumberFormat = "General"
$ws.Cells.Item(76, 5).Value = 1
$ws.Cells.Item(76, 6).NumberFormat = "General"
$ws.Cells.Item(76, 6).Value = 2
$ws.Cells.Item(76, 7).NumberFormat = "@"
$ws.Cells.Item(76, 7).Value = "Trend"
$ws.Cells.Item(76, 8).NumberFormat = "@"
$ws.Cells.Item(76, 8).Value = "16.95"
$ws.Cells.Item(76, 9).ClearContents()
$ws.Cells.Item(76, 10).ClearContents()
$ws.Cells.Item(76, 11).ClearContents()
$ws.Cells.Item(76, 12).ClearContents()
$ws.Cells.Item(76, 13).NumberFormat = "@"
$ws.Cells.Item(76, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Cells.Item(76, 14).NumberFormat = "@"
$ws.Cells.Item(76, 14).Value = "Trend USB-Stick 16 GB 16.95 Schweizer Franken"
$ws.Cells.Item(76, 15).NumberFormat = "@"
$ws.Cells.Item(76, 15).Value = "2022-07-28 20:59:25"

# Row 77: '6459240' - 'Trend USB-Stick 256 GB'
$ws.Cells.Item(77, 1).NumberFormat = "@"
$ws.Cells.Item(77, 1).Value = "6459240"
$ws.Cells.Item(77, 2).NumberFormat = "@"
$ws.Cells.Item(77, 2).Value = "Trend USB-Stick 256 GB"
$ws.Cells.Item(77, 3).NumberFormat = "@"
$ws.Cells.Item(77, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-256-gb/p/6459240"
$ws.Cells.Item(77, 4).ClearContents()
$ws.Cells.Item(77, 5).ClearContents()
$ws.Cells.Item(77, 6).NumberFormat = "General"
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).NumberFormat = "@"
$ws.Cells.Item(77, 7).Value = "Trend"
$ws.Cells.Item(77, 8).NumberFormat = "@"
$ws.Cells.Item(77, 8).Value = "59.95"
$ws.Cells.Item(77, 9).ClearContents()
$ws.Cells.Item(77, 10).ClearContents()
$ws.Cells.Item(77, 11).ClearContents()
$ws.Cells.Item(77, 12).ClearContents()
$ws.Cells.Item(77, 13).NumberFormat = "@"
$ws.Cells.Item(77, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Cells.Item(77, 14).NumberFormat = "@"
$ws.Cells.Item(77, 14).Value = "Trend USB-Stick 256 GB 59.95 Schweizer Franken"
$ws.Cells.Item(77, 15).NumberFormat = "@"
$ws.Cells.Item(77, 15).Value = "2022-07-28 20:59:25"

# Row 78: '6007538' - 'Trend USB-Stick 32 GB'
$ws.Cells.Item(78, 1).NumberFormat = "@"
$ws.Cells.Item(78, 1).Value = "6007538"
$ws.Cells.Item(78, 2).NumberFormat = "@"
$ws.Cells.Item(78, 2).Value = "Trend USB-Stick 32 GB"
$ws.Cells.Item(78, 3).NumberFormat = "@"
$ws.Cells.Item(78, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-32-gb/p/6007538"
$ws.Cells.Item(78, 4).ClearContents()
$ws.Cells.Item(78, 5).ClearContents()
$ws.Cells.Item(78, 6).NumberFormat = "General"
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).NumberFormat = "@"
$ws.Cells.Item(78, 7).Value = "Trend"
$ws.Cells.Item(78, 8).NumberFormat = "@"
$ws.Cells.Item(78, 8).Value = "29.95"
$ws.Cells.Item(78, 9).ClearContents()
$ws.Cells.Item(78, 10).ClearContents()
$ws.Cells.Item(78, 11).ClearContents()
$ws.Cells.Item(78, 12).ClearContents()
$ws.Cells.Item(78, 13).NumberFormat = "@"
$ws.Cells.Item(78, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Cells.Item(78, 14).NumberFormat = "@"
$ws.Cells.Item(78, 14).Value = "Trend USB-Stick 32 GB 29.95 Schweizer Franken"
$ws.Cells.Item(78, 15).NumberFormat = "@"
$ws.Cells.Item(78, 15).Value = "2022-07-28 20:59:25"

# Row 79: '5867973' - 'Trend USB-Stick 64 GB'
$ws.Cells.Item(79, 1).NumberFormat = "@"
$ws.Cells.Item(79, 1).Value = "5867973"
$ws.Cells.Item(79, 2).NumberFormat = "@"
$ws.Cells.Item(79, 2).Value = "Trend USB-Stick 64 GB"
$ws.Cells.Item(79, 3).NumberFormat = "@"
$ws.Cells.Item(79, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-64-gb/p/5867973"
$ws.Cells.Item(79, 4).ClearContents()
$ws.Cells.Item(79, 5).NumberFormat = "General"
$ws.Cells.Item(79, 5).Value = 4
$ws.Cells.Item(79, 6).NumberFormat = "General"
$ws.Cells.Item(79, 6).Value = 2
$ws.Cells.Item(79, 7).NumberFormat = "@"
$ws.Cells.Item(79, 7).Value = "Trend"
$ws.Cells.Item(79, 8).NumberFormat = "@"
$ws.Cells.Item(79, 8).Value = "39.95"
$ws.Cells.Item(79, 9).ClearContents()
$ws.Cells.Item(79, 10).ClearContents()
$ws.Cells.Item(79, 11).ClearContents()
$ws.Cells.Item(79, 12).ClearContents()
$ws.Cells.Item(79, 13).NumberFormat = "@"
$ws.Cells.Item(79, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Cells.Item(79, 14).NumberFormat = "@"
$ws.Cells.Item(79, 14).Value = "Trend USB-Stick 64 GB 39.95 Schweizer Franken"
$ws.Cells.Item(79, 15).NumberFormat = "@"
$ws.Cells.Item(79, 15).Value = "2022-07-28 20:59:25"

# Row 80: '5894674' - 'Trisa Beauty Sonic Nail Care System'
$ws.Cells.Item(80, 1).NumberFormat = "@"
$ws.Cells.Item(80, 1).Value = "5894674"
$ws.Cells.Item(80, 2).NumberFormat = "@"
$ws.Cells.Item(80, 2).Value = "Trisa Beauty Sonic Nail Care System"
$ws.Cells.Item(80, 3).NumberFormat = "@"
$ws.Cells.Item(80, 3).Value = "/de/kosmetik-gesundheit/make-up/naegel/set-accessoires/trisa-beauty-sonic-nail-care-system/p/5894674"
$ws.Cells.Item(80, 4).NumberFormat = "@"
$ws.Cells.Item(80, 4).Value = "1ST"
$ws.Cells.Item(80, 5).ClearContents()
$ws.Cells.Item(80, 6).NumberFormat = "General"
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).NumberFormat = "@"
$ws.Cells.Item(80, 7).Value = "Trisa"
$ws.Cells.Item(80, 8).NumberFormat = "@"
$ws.Cells.Item(80, 8).Value = "29.95"
$ws.Cells.Item(80, 9).NumberFormat = "@"
$ws.Cells.Item(80, 9).Value = "29.95/1ST"
$ws.Cells.Item(80, 10).NumberFormat = "@"
$ws.Cells.Item(80, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(80, 11).NumberFormat = "@"
$ws.Cells.Item(80, 11).Value = "29.95"
$ws.Cells.Item(80, 12).NumberFormat = "@"
$ws.Cells.Item(80, 12).Value = "1ST"
$ws.Cells.Item(80, 13).NumberFormat = "@"
$ws.Cells.Item(80, 13).Value = "['kosmetik-gesundheit', 'make-up', 'naegel', 'set-accessoires']"
$ws.Cells.Item(80, 14).NumberFormat = "@"
$ws.Cells.Item(80, 14).Value = "Trisa Beauty Sonic Nail Care System 29.95 Schweizer Franken"
$ws.Cells.Item(80, 15).NumberFormat = "@"
$ws.Cells.Item(80, 15).Value = "2022-07-28 20:59:25"

# Row 81: '4589934' - 'Varta Longlife AA 4er Bli'
$ws.Cells.Item(81, 1).NumberFormat = "@"
$ws.Cells.Item(81, 1).Value = "4589934"
$ws.Cells.Item(81, 2).NumberFormat = "@"
$ws.Cells.Item(81, 2).Value = "Varta Longlife AA 4er Bli"
$ws.Cells.Item(81, 3).NumberFormat = "@"
$ws.Cells.Item(81, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-aa-4er-bli/p/4589934"
$ws.Cells.Item(81, 4).NumberFormat = "@"
$ws.Cells.Item(81, 4).Value = "4ST"
$ws.Cells.Item(81, 5).ClearContents()
$ws.Cells.Item(81, 6).NumberFormat = "General"
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).NumberFormat = "@"
$ws.Cells.Item(81, 7).Value = "Varta"
$ws.Cells.Item(81, 8).NumberFormat = "@"
$ws.Cells.Item(81, 8).Value = "7.95"
$ws.Cells.Item(81, 9).NumberFormat = "@"
$ws.Cells.Item(81, 9).Value = "1.99/1ST"
$ws.Cells.Item(81, 10).NumberFormat = "@"
$ws.Cells.Item(81, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(81, 11).NumberFormat = "@"
$ws.Cells.Item(81, 11).Value = "1.99"
$ws.Cells.Item(81, 12).NumberFormat = "@"
$ws.Cells.Item(81, 12).Value = "1ST"
$ws.Cells.Item(81, 13).NumberFormat = "@"
$ws.Cells.Item(81, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Cells.Item(81, 14).NumberFormat = "@"
$ws.Cells.Item(81, 14).Value = "Varta Longlife AA 4er Bli 7.95 Schweizer Franken"
$ws.Cells.Item(81, 15).NumberFormat = "@"
$ws.Cells.Item(81, 15).Value = "2022-07-28 20:59:25"

# Row 82: '4589933' - 'Varta Longlife AAA 4er Bli'
$ws.Cells.Item(82, 1).NumberFormat = "@"
$ws.Cells.Item(82, 1).Value = "4589933"
$ws.Cells.Item(82, 2).NumberFormat = "@"
$ws.Cells.Item(82, 2).Value = "Varta Longlife AAA 4er Bli"
$ws.Cells.Item(82, 3).NumberFormat = "@"
$ws.Cells.Item(82, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-aaa-4er-bli/p/4589933"
$ws.Cells.Item(82, 4).NumberFormat = "@"
$ws.Cells.Item(82, 4).Value = "4ST"
$ws.Cells.Item(82, 5).ClearContents()
$ws.Cells.Item(82, 6).NumberFormat = "General"
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).NumberFormat = "@"
$ws.Cells.Item(82, 7).Value = "Varta"
$ws.Cells.Item(82, 8).NumberFormat = "@"
$ws.Cells.Item(82, 8).Value = "7.95"
$ws.Cells.Item(82, 9).NumberFormat = "@"
$ws.Cells.Item(82, 9).Value = "1.99/1ST"
$ws.Cells.Item(82, 10).NumberFormat = "@"
$ws.Cells.Item(82, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(82, 11).NumberFormat = "@"
$ws.Cells.Item(82, 11).Value = "1.99"
$ws.Cells.Item(82, 12).NumberFormat = "@"
$ws.Cells.Item(82, 12).Value = "1ST"
$ws.Cells.Item(82, 13).NumberFormat = "@"
$ws.Cells.Item(82, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Cells.Item(82, 14).NumberFormat = "@"
$ws.Cells.Item(82, 14).Value = "Varta Longlife AAA 4er Bli 7.95 Schweizer Franken"
$ws.Cells.Item(82, 15).NumberFormat = "@"
$ws.Cells.Item(82, 15).Value = "2022-07-28 20:59:25"

# Row 83: '4589935' - 'Varta Longlife C 2er Bli'
$ws.Cells.Item(83, 1).NumberFormat = "@"
$ws.Cells.Item(83, 1).Value = "4589935"
$ws.Cells.Item(83, 2).NumberFormat = "@"
$ws.Cells.Item(83, 2).Value = "Varta Longlife C 2er Bli"
$ws.Cells.Item(83, 3).NumberFormat = "@"
$ws.Cells.Item(83, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-c-2er-bli/p/4589935"
$ws.Cells.Item(83, 4).NumberFormat = "@"
$ws.Cells.Item(83, 4).Value = "2ST"
$ws.Cells.Item(83, 5).ClearContents()
$ws.Cells.Item(83, 6).NumberFormat = "General"
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).NumberFormat = "@"
$ws.Cells.Item(83, 7).Value = "Varta"
$ws.Cells.Item(83, 8).NumberFormat = "@"
$ws.Cells.Item(83, 8).Value = "6.95"
$ws.Cells.Item(83, 9).NumberFormat = "@"
$ws.Cells.Item(83, 9).Value = "3.48/1ST"
$ws.Cells.Item(83, 10).NumberFormat = "@"
$ws.Cells.Item(83, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(83, 11).NumberFormat = "@"
$ws.Cells.Item(83, 11).Value = "3.48"
$ws.Cells.Item(83, 12).NumberFormat = "@"
$ws.Cells.Item(83, 12).Value = "1ST"
$ws.Cells.Item(83, 13).NumberFormat = "@"
$ws.Cells.Item(83, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(83, 14).NumberFormat = "@"
$ws.Cells.Item(83, 14).Value = "Varta Longlife C 2er Bli 6.95 Schweizer Franken"
$ws.Cells.Item(83, 15).NumberFormat = "@"
$ws.Cells.Item(83, 15).Value = "2022-07-28 20:59:25"

# Row 84: '3494067' - 'Varta Longlife Max Power AAA 4er Bli'
$ws.Cells.Item(84, 1).NumberFormat = "@"
$ws.Cells.Item(84, 1).Value = "3494067"
$ws.Cells.Item(84, 2).NumberFormat = "@"
$ws.Cells.Item(84, 2).Value = "Varta Longlife Max Power AAA 4er Bli"
$ws.Cells.Item(84, 3).NumberFormat = "@"
$ws.Cells.Item(84, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-max-power-aaa-4er-bli/p/3494067"
$ws.Cells.Item(84, 4).NumberFormat = "@"
$ws.Cells.Item(84, 4).Value = "4ST"
$ws.Cells.Item(84, 5).ClearContents()
$ws.Cells.Item(84, 6).NumberFormat = "General"
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).NumberFormat = "@"
$ws.Cells.Item(84, 7).Value = "Varta"
$ws.Cells.Item(84, 8).NumberFormat = "@"
$ws.Cells.Item(84, 8).Value = "9.95"
$ws.Cells.Item(84, 9).NumberFormat = "@"
$ws.Cells.Item(84, 9).Value = "2.49/1ST"
$ws.Cells.Item(84, 10).NumberFormat = "@"
$ws.Cells.Item(84, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(84, 11).NumberFormat = "@"
$ws.Cells.Item(84, 11).Value = "2.49"
$ws.Cells.Item(84, 12).NumberFormat = "@"
$ws.Cells.Item(84, 12).Value = "1ST"
$ws.Cells.Item(84, 13).NumberFormat = "@"
$ws.Cells.Item(84, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Cells.Item(84, 14).NumberFormat = "@"
$ws.Cells.Item(84, 14).Value = "Varta Longlife Max Power AAA 4er Bli 9.95 Schweizer Franken"
$ws.Cells.Item(84, 15).NumberFormat = "@"
$ws.Cells.Item(84, 15).Value = "2022-07-28 20:59:25"

# Row 85: '3591269' - 'Varta Longlife Max Power C 2er Bli'
$ws.Cells.Item(85, 1).NumberFormat = "@"
$ws.Cells.Item(85, 1).Value = "3591269"
$ws.Cells.Item(85, 2).NumberFormat = "@"
$ws.Cells.Item(85, 2).Value = "Varta Longlife Max Power C 2er Bli"
$ws.Cells.Item(85, 3).NumberFormat = "@"
$ws.Cells.Item(85, 3).Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-max-power-c-2er-bli/p/3591269"
$ws.Cells.Item(85, 4).NumberFormat = "@"
$ws.Cells.Item(85, 4).Value = "2ST"
$ws.Cells.Item(85, 5).NumberFormat = "General"
$ws.Cells.Item(85, 5).Value = 1
$ws.Cells.Item(85, 6).NumberFormat = "General"
$ws.Cells.Item(85, 6).Value = 5
$ws.Cells.Item(85, 7).NumberFormat = "@"
$ws.Cells.Item(85, 7).Value = "Varta"
$ws.Cells.Item(85, 8).NumberFormat = "@"
$ws.Cells.Item(85, 8).Value = "8.95"
$ws.Cells.Item(85, 9).NumberFormat = "@"
$ws.Cells.Item(85, 9).Value = "4.48/1ST"
$ws.Cells.Item(85, 10).NumberFormat = "@"
$ws.Cells.Item(85, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(85, 11).NumberFormat = "@"
$ws.Cells.Item(85, 11).Value = "4.48"
$ws.Cells.Item(85, 12).NumberFormat = "@"
$ws.Cells.Item(85, 12).Value = "1ST"
$ws.Cells.Item(85, 13).NumberFormat = "@"
$ws.Cells.Item(85, 13).Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Cells.Item(85, 14).NumberFormat = "@"
$ws.Cells.Item(85, 14).Value = "Varta Longlife Max Power C 2er Bli 8.95 Schweizer Franken"
$ws.Cells.Item(85, 15).NumberFormat = "@"
$ws.Cells.Item(85, 15).Value = "2022-07-28 20:59:25"

